# The deck's slide master currently uses the "Integral" theme
# (ppt/theme/theme1.xml). This commit swaps the presentation over to the
# stock "Office Theme" palette (the colours that already live, unused, in
# ppt/theme/theme2.xml, which only the notes master points at).
#
# PowerPoint's theme-color API is exposed on the slide object as
# ThemeColorScheme - a 12-entry collection in the fixed clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# Writing .RGB for each entry rewrites <a:clrScheme> in the theme part
# backing the slide master, which is exactly the edit the diff shows.

function Convert-HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB long is little-endian 0x00BBGGRR.
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" colours, in clrScheme order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Convert-HexToRgb $officeThemeColors[$i - 1]
}
